$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. These cells hold numeric-looking figures that are
# stored as TEXT (shared strings) in the workbook, so we force text entry
# (leading apostrophe) and restore the cell's original style afterwards -
# otherwise Excel would either (a) reinterpret the numeric-looking string as
# a real number, or (b) tag the cell with an extra "quote prefix" style.
$updates = @{
    "D11" = "9.36"
    "B33" = "6.96"
    "D33" = "8.47"
    "B36" = "77.48"
    "C36" = "16.74"
    "D36" = "94.21"
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.Value = "'" + $updates[$addr]
    $rng.Style = $origStyle
}
